$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix SOCP constraint handling: reset Bs (H) and Gs (I) columns to 0 for rows 8-16
$ws.Range("H8:I16").Value = 0

# Update the active selection on the sheet to I22 (was H22)
$ws.Range("I22").Select()

$wb.Save()
